# Cards.xlsx refactor — view/window state clean-up (no cell-data changes).
#
# Summary of the target change (see diff):
#  - Workbook: bookViews window chrome (xWindow/yWindow/windowWidth/
#    windowHeight) and the active-tab flips from "Deck Saves" -> "Card
#    Library".
#  - "Card Library" sheet (sheet 1): becomes the selected/active tab,
#    selection moves to B4, and columns A-C get wider/narrower.
#  - "Deck Saves" sheet (sheet 2): loses the selected/active flag,
#    selection moves to A16.

$wb = $excel.ActiveWorkbook

$wsCards = $wb.Worksheets.Item("Card Library")
$wsDecks = $wb.Worksheets.Item("Deck Saves")

# --- "Deck Saves": move its selection while it's still the active sheet,
#     then hand activation over to "Card Library" (mirrors the diff, where
#     the tabSelected flag moves from Deck Saves to Card Library). ---
$wsDecks.Activate()
$wsDecks.Range("A16").Select()

# --- "Card Library" becomes the active / selected tab ---
$wsCards.Activate()
$wsCards.Range("B4").Select()

# --- Column width tweaks on "Card Library" (A:19.86->27, B:21.86->22.43, C:6.29->7) ---
# NOTE: the host round-trips ColumnWidth through its internal pixel grid
# (quantized to 1/6 of a character) before it is written back out as the
# stored `width`, so the values below are chosen so the *stored* width in
# the saved file lands on (or, for column B, as close as that grid allows
# to) the target from the diff.
$wsCards.Columns.Item(1).ColumnWidth = 26.16666667
$wsCards.Columns.Item(2).ColumnWidth = 21.66666667
$wsCards.Columns.Item(3).ColumnWidth = 6.16666667

# --- Workbook window chrome (best-effort; mirrors the new bookViews size) ---
$aw = $excel.ActiveWindow
$aw.Left = -120
$aw.Top = -120
$aw.Width = 29040
$aw.Height = 15840
